# Automatische test-sync: 2025-06-24 21:06:50
# Adds the new "Productmaat ruilen" mail-log entry as row 27 on the "Logs"
# sheet, extends the conditional formatting ranges to include the new row,
# and refreshes the "Retour / Terugbetaling" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A27").Value = "Productmaat ruilen"
$ws.Range("B27").Value = "mailmind.test@zohomail.eu"
$ws.Range("C27").Value = "Wil graag andere maat ontvangen`nSent using {0}"
$ws.Range("D27").Value = "Retour / Terugbetaling"
$ws.Range("E27").Value = "Beste klant,`nBedankt voor je bericht. Om je beter van dienst te kunnen zijn, vragen we je om wat meer informatie te verstrekken. Zoals de bestelgegevens, zoals het ordernummer en de huidige maat van het product. Op deze manier kunnen we je verzoek zo snel mogelijk verwerken. `nWe zien je reactie graag tegemoet.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$ws.Range("F27").Value = "2025-06-24 21:06:06"
$ws.Range("G27").Value = "Ja"

# Extend the two conditional-formatting ranges from row 26 to row 27.
$ws.Range("D2:D26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D27"))
$ws.Range("G2:G26").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G27"))

# Update the Dashboard summary count for "Retour / Terugbetaling" (9 -> 10).
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 10
